$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.280.88'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.29%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.791.65'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.62%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.54'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.80%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.594'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.97%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '35.97'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.65%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.289'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.35%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0671'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.04%  '

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.87%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.053.76'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.41%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.09'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.94%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.814.33'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.15%  '

# Row 15
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '34.308.13'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.06%  '

# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.625'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.45%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.35'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.22%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.56'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.03%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '239.58'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.98%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0765'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.50%  '

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.03%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.18'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.99%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.05'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.33%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.18'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.28%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '170.31'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.10%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.90'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +4.23%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.11'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.36%  '

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.59%  '

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.02%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.22'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.89%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.75'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.59%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.85'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.43%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0509'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.90%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.74'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -5.55%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.354.39'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.39%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.640'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.94%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.04'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.57%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.32'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -9.82%  '

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.15%  '

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.31%  '

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.62%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '80.21'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.48%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.925'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.00%  '

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +5.42%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.08'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -5.81%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0496'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.81%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.955.90'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.38%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.75'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.03%  '

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.09%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '101.27'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.84%  '

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0601'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.57%  '
